$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row of data: "python_project_structure_with_explanation" with count 5
$ws.Range("A22").Value = "python_project_structure_with_explanation"
$ws.Range("B22").Value = 5

# Update the active selection to B23 (as in the final sheet view)
$ws.Range("B23").Select()
